# Update "想去人数" (interested-count) figures to the values captured at the
# later data-refresh snapshot (commit 456a3b4).
#
# Sheet "展览"   (Exhibitions)      -> column F on rows 5,7,8,9,10,13
# Sheet "本地生活" (Local life)     -> column F on rows 2,3
# Sheet "全部类型" (All categories) -> column F on rows 2,3,13,17,19,20,22,27

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 789
$wsExhibition.Range("F7").Value = 528
$wsExhibition.Range("F8").Value = 101
$wsExhibition.Range("F9").Value = 547
$wsExhibition.Range("F10").Value = 492
$wsExhibition.Range("F13").Value = 134

$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Range("F2").Value = 6246
$wsLocalLife.Range("F3").Value = 767

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 6246
$wsAllTypes.Range("F3").Value = 767
$wsAllTypes.Range("F13").Value = 789
$wsAllTypes.Range("F17").Value = 528
$wsAllTypes.Range("F19").Value = 101
$wsAllTypes.Range("F20").Value = 547
$wsAllTypes.Range("F22").Value = 492
$wsAllTypes.Range("F27").Value = 134
